{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Delete the trailing link-attribution paragraphs (old \"\u05de\u05d0\u05de\u05e8 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d5\u05de\u05d5\u05de\u05dc\u05e5\n// \u05dc\u05e7\u05e8\u05d9\u05d0\u05d4!\" and the old arxiv link paragraph) - they are removed entirely in\n// the revised review. Delete from the end so earlier indices stay valid.\nparagraphs.items[9].delete();\nparagraphs.items[8].delete();\nawait context.sync();\n\n// Re-load the remaining paragraphs so we have fresh, stable references.\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Each target run is written with clear() + insertText(\"Replace\") rather\n// than insertText() alone: some of the original paragraphs carried\n// xml:space=\"preserve\" (trailing-space runs) and a bare Replace keeps that\n// stale attribute even once the new text no longer needs it. clear()\n// first drops the old run (and its whitespace flag) so the freshly\n// inserted text gets space-preservation recomputed from its own content.\n\n// Paragraph 0: update the date in the daily-review header.\nparagraphs.items[0].clear();\nparagraphs.items[0].insertText(\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 01.09.24: \u26a1\ufe0f\ud83d\ude80\",\n  \"Replace\"\n);\n\n// Paragraph 1: replace the paper title (also drops the trailing line break\n// that followed the old title run).\nparagraphs.items[1].clear();\nparagraphs.items[1].insertText(\n  \"DIFFUSION MODELS ARE REAL-TIME GAME ENGINES\",\n  \"Replace\"\n);\n\n// Paragraph 2: intro paragraph rewritten for the new paper.\nparagraphs.items[2].clear();\nparagraphs.items[2].insertText(\n  \"\u05d8\u05d5\u05d1, \u05e2\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05e4\u05e9\u05d5\u05d8 \u05dc\u05d0 \u05d4\u05d9\u05d4 \u05dc\u05d3\u05dc\u05d2 \u05de\u05db\u05de\u05d4 \u05e1\u05d9\u05d1\u05d5\u05ea. \u05d4\u05e1\u05d9\u05d1\u05d4 \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 \u05e9\u05d0\u05e0\u05d9 \u05de\u05e1\u05e4\u05d9\u05e7 \u05e2\u05ea\u05d9\u05e7 \u05d5\u05e2\u05d5\u05d3 \u05e9\u05d9\u05d7\u05e7\u05ea\u05d9 \u05d1\u05de\u05e9\u05d7\u05e7 \u05d4\u05e0\u05e7\u05e8\u05d0 \u05d3\u05d5\u05dd (doom) \u05d1\u05de\u05d5 \u05d9\u05d3\u05d9\u05d9 \u05db\u05d0\u05e9\u05e8 \u05d4\u05d9\u05d9\u05ea\u05d9 \u05e0\u05e2\u05e8. \u05d3\u05d1\u05e8 \u05e9\u05e0\u05d9 \u05dc\u05d0 \u05db\u05dc \u05d9\u05d5\u05dd \u05de\u05d7\u05dc\u05d9\u05e4\u05d9\u05dd \u05dc\u05da \u05de\u05e0\u05d5\u05e2 \u05de\u05e9\u05d7\u05e7 \u05d1\u05de\u05d5\u05d3\u05dc \u05dc\u05de\u05d9\u05d3\u05ea \u05de\u05db\u05d5\u05e0\u05ea \u05d0\u05d5 \u05d1\u05e9\u05de\u05d5 \u05d4\u05de\u05d5\u05db\u05e8 AI. \u05db\u05de\u05d5\u05d1\u05df \u05e9\u05d6\u05d4 \u05db\u05d9\u05d5\u05d5\u05df \u05de\u05d7\u05e7\u05e8 \u05de\u05d0\u05d5\u05d3 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e2\u05dd \u05e4\u05d5\u05d8\u05e0\u05e6\u05d9\u05d0\u05dc \u05dc\u05d4\u05ea\u05e4\u05ea\u05d7 \u05dc\u05db\u05dc\u05d9\u05dd \u05de\u05d1\u05d5\u05e1\u05e1\u05d9 AI \u05dc\u05d1\u05e0\u05d9\u05d9\u05ea \u05de\u05e9\u05d7\u05e7\u05d9 \u05de\u05d7\u05e9\u05d1 \u05d7\u05d3\u05e9\u05d9\u05dd.\",\n  \"Replace\"\n);\n\n// Paragraph 3.\nparagraphs.items[3].clear();\nparagraphs.items[3].insertText(\n  \"\u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05e0\u05d5 \u05d3\u05d9 \u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d9\u05d8\u05d9\u05d1\u05d9. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05e1\u05d5\u05db\u05df (agent) \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05e9\u05d7\u05e7 \u05de\u05e9\u05d7\u05e7 \u05d3\u05d5\u05dd \u05d1\u05e2\u05e6\u05de\u05d5 \u05e2\u05dc \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05d4\u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05e9\u05e9\u05d5\u05d7\u05e7\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d1\u05e0\u05d9 \u05d0\u05d3\u05dd. \u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05db\u05de\u05d4 \u05de\u05de\u05e6\u05d1\u05d9 \u05d4\u05de\u05e9\u05d7\u05e7 (\u05e4\u05e8\u05d9\u05d9\u05de\u05d9\u05dd) \u05d5\u05d4\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d5\u05ea (\u05d9\u05e8\u05d9, \u05ea\u05e0\u05d5\u05e2\u05d4, \u05e4\u05d2\u05d9\u05e2\u05d4 \u05d5\u05db\u05d3\u05d5\u05de\u05d4) \u05de\u05d8\u05e8\u05ea \u05d4\u05e1\u05d5\u05db\u05df \u05d4\u05d9\u05d0 \u05d7\u05d9\u05d6\u05d5\u05d9 \u05d4\u05e4\u05e2\u05d5\u05dc\u05ea\u05d5 \u05d4\u05d1\u05d0\u05d4. \u05d6\u05d4 \u05e0\u05e2\u05e9\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d8\u05db\u05e0\u05d9\u05e7\u05d5\u05ea RL \u05d3\u05d9 \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05d5\u05ea \u05db\u05d0\u05e9\u05e8 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4-reward \u05e0\u05d1\u05d7\u05e8\u05d4 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d4\u05d2\u05d9\u05d5\u05e0\u05d9\u05ea \u05d1\u05d4\u05ea\u05d0\u05dd \u05dc\u05dc\u05d5\u05d2\u05d9\u05e7\u05ea \u05d4\u05de\u05e9\u05d7\u05e7 (\u05db\u05dc\u05d5\u05de\u05e8 \u05e4\u05d2\u05d9\u05e2\u05d4 \u05d0\u05d5 \u05de\u05d5\u05d5\u05ea \u05e9\u05dc \u05d4\u05e1\u05d5\u05db\u05df \u05de\u05e7\u05d1\u05dc\u05d5\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05e9\u05dc\u05d9\u05e9\u05d9 \u05d5\u05d0\u05d9\u05dc\u05d5 \u05e4\u05d2\u05d9\u05e2\u05d4 \u05d1\u05d0\u05d5\u05d9\u05d1, \u05d0\u05d9\u05e1\u05d5\u05e3 \u05e0\u05e9\u05e7 \u05d5\u05db\u05d3\u05d5\u05de\u05d4 \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05ea\u05d2\u05de\u05d5\u05dc \u05d7\u05d9\u05d5\u05d1\u05d9).\",\n  \"Replace\"\n);\n\n// Paragraph 4 (keeps a trailing space, like the original run).\nparagraphs.items[4].clear();\nparagraphs.items[4].insertText(\n  \"\u05d0\u05d7\u05e8\u05d9 \u05e9\u05d4\u05e1\u05d5\u05db\u05df \u05dc\u05de\u05d3 \u05dc\u05e9\u05d7\u05e7 \u05d3\u05d5\u05dd, \u05de\u05d2\u05e0\u05e8\u05d8\u05d9\u05dd \u05db\u05de\u05d5\u05ea \u05de\u05d0\u05d5\u05d3 \u05d2\u05d3\u05d5\u05dc\u05d4 \u05e9\u05dc \u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05d3\u05d5\u05dd \u05e2\u05dd \u05d4\u05e1\u05d5\u05db\u05df. \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05e1\u05d5\u05db\u05df \u05de\u05e9\u05d7\u05e7 \u05d1\u05de\u05e9\u05d7\u05e7 \u05d0\u05de\u05d9\u05ea\u05d9 \u05db\u05de\u05d5 \u05d0\u05d7\u05d3 \u05d4\u05d0\u05d3\u05dd. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05e8\u05d9\u05d9\u05dd \u05d4\u05d1\u05d0 \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05e4\u05e8\u05d9\u05d9\u05de\u05d9\u05dd \u05d4\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d4\u05e7\u05d5\u05d3\u05de\u05d5\u05ea \u05d5\u05d4\u05e0\u05d5\u05db\u05d7\u05d9\u05ea. \",\n  \"Replace\"\n);\n\n// Paragraph 5.\nparagraphs.items[5].clear();\nparagraphs.items[5].insertText(\n  \"\u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05ea\u05d1\u05e6\u05e2 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d9 \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea: \u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05de\u05e7\u05d1\u05dc \u05db\u05e7\u05dc\u05d8 \u05d0\u05ea \u05d4\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d4\u05e7\u05d5\u05d3\u05de\u05d5\u05ea \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d3\u05d5\u05e8 (\u05e9\u05de\u05d0\u05d5\u05de\u05df \u05d2\u05dd \u05db\u05df) \u05d5\u05d1\u05e0\u05d5\u05e1\u05e3 \u05d0\u05ea \u05d4\u05e4\u05e8\u05d9\u05d9\u05de\u05d9\u05dd \u05d4\u05e7\u05d5\u05d3\u05de\u05d9\u05dd \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc\u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 (\u05d1\u05e6\u05d5\u05e8\u05d4 \u05de\u05d5\u05e8\u05e2\u05e9\u05ea \u05dc\u05e9\u05d9\u05e4\u05d5\u05e8 \u05d9\u05db\u05d5\u05dc\u05ea \u05d4\u05db\u05dc\u05dc\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc). \u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05e9\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d4\u05e9\u05ea\u05de\u05e9\u05d5 \u05d1\u05d5 \u05d4\u05d9\u05e0\u05d5 \u05dc\u05d8\u05e0\u05d8\u05d9 (\u05db\u05dc\u05d5\u05de\u05e8 \u05d7\u05d9\u05d6\u05d5\u05d9 \u05d4\u05e8\u05e2\u05e9 \u05de\u05ea\u05d1\u05e6\u05e2 \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc \u05d4\u05e4\u05e8\u05d9\u05d9\u05dd \u05d4\u05e0\u05d7\u05d6\u05d4). \u05e0\u05e6\u05d9\u05d9\u05df \u05db\u05d9 \u05db\u05d0\u05df \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de\u05de\u05d5\u05d3\u05dc\u05d9 \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d9\u05e9\u05e0\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d5\u05d3\u05dc \u05d4\u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d1\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05de\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0 \u05f4\u05de\u05d4\u05d9\u05e8\u05d5\u05ea\u05f4 \u05e9\u05dc \u05d4\u05e4\u05e8\u05d9\u05d9\u05dd \u05d4\u05de\u05d5\u05e8\u05e2\u05e9 \u05e9\u05d4\u05d9\u05d0 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d4\u05e4\u05e8\u05d9\u05d9\u05dd \u05d4\u05e0\u05e7\u05d9 \u05d5\u05d4\u05e8\u05e2\u05e9 \u05d4\u05de\u05ea\u05d5\u05d5\u05e1\u05e3 \u05d0\u05dc\u05d9\u05d5 \u05d1\u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d4. \u05e8\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d6\u05d5 \u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d0\u05d9\u05db\u05d5\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05d5\u05de\u05d0\u05d9\u05e6\u05d4 \u05d4\u05ea\u05db\u05e0\u05e1\u05d5\u05ea\u05d4 (\u05de\u05d5\u05db\u05d7 \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05ea \u05db\u05e8\u05d2\u05d9\u05dc)...\",\n  \"Replace\"\n);\n\n// Paragraph 6: short closing remark.\nparagraphs.items[6].clear();\nparagraphs.items[6].insertText(\"\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05d5\u05d3 \u05de\u05d2\u05e0\u05d9\u05d1\u2026\", \"Replace\");\n\n// Paragraph 7: new arxiv link (trailing space preserved, like the source).\nparagraphs.items[7].clear();\nparagraphs.items[7].insertText(\n  \"https://arxiv.org/pdf/2408.14837 \",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Drop the two trailing paragraphs entirely (old \"\u05de\u05d0\u05de\u05e8 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05d5\u05de\u05d5\u05de\u05dc\u05e5\n# \u05dc\u05e7\u05e8\u05d9\u05d0\u05d4!\" line and the old arxiv link line) - the revised review folds\n# the link into the body text above and removes both. Delete from the\n# end so earlier paragraph indices stay valid.\n$d.Paragraphs.Item(10).Range.Delete()\n$d.Paragraphs.Item(9).Range.Delete()\n\n# For every remaining paragraph, replace the run text but exclude the\n# trailing paragraph-mark character from the range first (MoveEnd -1)\n# before assigning .Text. That avoids merging with the following\n# paragraph (which a Delete()+InsertAfter() on the full paragraph range\n# would do) and - because the run is produced fresh rather than reusing\n# the old run's whitespace flag - xml:space=\"preserve\" only shows up on\n# the output when the new text actually has leading/trailing whitespace.\n\n$r1 = $d.Paragraphs.Item(1).Range\n[void]$r1.MoveEnd(1, -1)\n$r1.Text = \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 01.09.24: \u26a1\ufe0f\ud83d\ude80\"\n\n$r2 = $d.Paragraphs.Item(2).Range\n[void]$r2.MoveEnd(1, -1)\n$r2.Text = \"DIFFUSION MODELS ARE REAL-TIME GAME ENGINES\"\n\n$r3 = $d.Paragraphs.Item(3).Range\n[void]$r3.MoveEnd(1, -1)\n$r3.Text = \"\u05d8\u05d5\u05d1, \u05e2\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d6\u05d4 \u05e4\u05e9\u05d5\u05d8 \u05dc\u05d0 \u05d4\u05d9\u05d4 \u05dc\u05d3\u05dc\u05d2 \u05de\u05db\u05de\u05d4 \u05e1\u05d9\u05d1\u05d5\u05ea. \u05d4\u05e1\u05d9\u05d1\u05d4 \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 \u05e9\u05d0\u05e0\u05d9 \u05de\u05e1\u05e4\u05d9\u05e7 \u05e2\u05ea\u05d9\u05e7 \u05d5\u05e2\u05d5\u05d3 \u05e9\u05d9\u05d7\u05e7\u05ea\u05d9 \u05d1\u05de\u05e9\u05d7\u05e7 \u05d4\u05e0\u05e7\u05e8\u05d0 \u05d3\u05d5\u05dd (doom) \u05d1\u05de\u05d5 \u05d9\u05d3\u05d9\u05d9 \u05db\u05d0\u05e9\u05e8 \u05d4\u05d9\u05d9\u05ea\u05d9 \u05e0\u05e2\u05e8. \u05d3\u05d1\u05e8 \u05e9\u05e0\u05d9 \u05dc\u05d0 \u05db\u05dc \u05d9\u05d5\u05dd \u05de\u05d7\u05dc\u05d9\u05e4\u05d9\u05dd \u05dc\u05da \u05de\u05e0\u05d5\u05e2 \u05de\u05e9\u05d7\u05e7 \u05d1\u05de\u05d5\u05d3\u05dc \u05dc\u05de\u05d9\u05d3\u05ea \u05de\u05db\u05d5\u05e0\u05ea \u05d0\u05d5 \u05d1\u05e9\u05de\u05d5 \u05d4\u05de\u05d5\u05db\u05e8 AI. \u05db\u05de\u05d5\u05d1\u05df \u05e9\u05d6\u05d4 \u05db\u05d9\u05d5\u05d5\u05df \u05de\u05d7\u05e7\u05e8 \u05de\u05d0\u05d5\u05d3 \u05de\u05e2\u05e0\u05d9\u05d9\u05df \u05e2\u05dd \u05e4\u05d5\u05d8\u05e0\u05e6\u05d9\u05d0\u05dc \u05dc\u05d4\u05ea\u05e4\u05ea\u05d7 \u05dc\u05db\u05dc\u05d9\u05dd \u05de\u05d1\u05d5\u05e1\u05e1\u05d9 AI \u05dc\u05d1\u05e0\u05d9\u05d9\u05ea \u05de\u05e9\u05d7\u05e7\u05d9 \u05de\u05d7\u05e9\u05d1 \u05d7\u05d3\u05e9\u05d9\u05dd.\"\n\n$r4 = $d.Paragraphs.Item(4).Range\n[void]$r4.MoveEnd(1, -1)\n$r4.Text = \"\u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05e9\u05dc \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05e0\u05d5 \u05d3\u05d9 \u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d9\u05d8\u05d9\u05d1\u05d9. \u05d1\u05e9\u05dc\u05d1 \u05d4\u05e8\u05d0\u05e9\u05d5\u05df \u05d4\u05e1\u05d5\u05db\u05df (agent) \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05e9\u05d7\u05e7 \u05de\u05e9\u05d7\u05e7 \u05d3\u05d5\u05dd \u05d1\u05e2\u05e6\u05de\u05d5 \u05e2\u05dc \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05d4\u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05e9\u05e9\u05d5\u05d7\u05e7\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d1\u05e0\u05d9 \u05d0\u05d3\u05dd. \u05db\u05dc\u05d5\u05de\u05e8 \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05db\u05de\u05d4 \u05de\u05de\u05e6\u05d1\u05d9 \u05d4\u05de\u05e9\u05d7\u05e7 (\u05e4\u05e8\u05d9\u05d9\u05de\u05d9\u05dd) \u05d5\u05d4\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05d7\u05e8\u05d5\u05e0\u05d5\u05ea (\u05d9\u05e8\u05d9, \u05ea\u05e0\u05d5\u05e2\u05d4, \u05e4\u05d2\u05d9\u05e2\u05d4 \u05d5\u05db\u05d3\u05d5\u05de\u05d4) \u05de\u05d8\u05e8\u05ea \u05d4\u05e1\u05d5\u05db\u05df \u05d4\u05d9\u05d0 \u05d7\u05d9\u05d6\u05d5\u05d9 \u05d4\u05e4\u05e2\u05d5\u05dc\u05ea\u05d5 \u05d4\u05d1\u05d0\u05d4. \u05d6\u05d4 \u05e0\u05e2\u05e9\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05d8\u05db\u05e0\u05d9\u05e7\u05d5\u05ea RL \u05d3\u05d9 \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05d5\u05ea \u05db\u05d0\u05e9\u05e8 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4-reward \u05e0\u05d1\u05d7\u05e8\u05d4 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d4\u05d2\u05d9\u05d5\u05e0\u05d9\u05ea \u05d1\u05d4\u05ea\u05d0\u05dd \u05dc\u05dc\u05d5\u05d2\u05d9\u05e7\u05ea \u05d4\u05de\u05e9\u05d7\u05e7 (\u05db\u05dc\u05d5\u05de\u05e8 \u05e4\u05d2\u05d9\u05e2\u05d4 \u05d0\u05d5 \u05de\u05d5\u05d5\u05ea \u05e9\u05dc \u05d4\u05e1\u05d5\u05db\u05df \u05de\u05e7\u05d1\u05dc\u05d5\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05e9\u05dc\u05d9\u05e9\u05d9 \u05d5\u05d0\u05d9\u05dc\u05d5 \u05e4\u05d2\u05d9\u05e2\u05d4 \u05d1\u05d0\u05d5\u05d9\u05d1, \u05d0\u05d9\u05e1\u05d5\u05e3 \u05e0\u05e9\u05e7 \u05d5\u05db\u05d3\u05d5\u05de\u05d4 \u05de\u05e7\u05d1\u05dc\u05d9\u05dd \u05ea\u05d2\u05de\u05d5\u05dc \u05d7\u05d9\u05d5\u05d1\u05d9).\"\n\n$r5 = $d.Paragraphs.Item(5).Range\n[void]$r5.MoveEnd(1, -1)\n$r5.Text = \"\u05d0\u05d7\u05e8\u05d9 \u05e9\u05d4\u05e1\u05d5\u05db\u05df \u05dc\u05de\u05d3 \u05dc\u05e9\u05d7\u05e7 \u05d3\u05d5\u05dd, \u05de\u05d2\u05e0\u05e8\u05d8\u05d9\u05dd \u05db\u05de\u05d5\u05ea \u05de\u05d0\u05d5\u05d3 \u05d2\u05d3\u05d5\u05dc\u05d4 \u05e9\u05dc \u05de\u05e9\u05d7\u05e7\u05d9\u05dd \u05d3\u05d5\u05dd \u05e2\u05dd \u05d4\u05e1\u05d5\u05db\u05df. \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05e1\u05d5\u05db\u05df \u05de\u05e9\u05d7\u05e7 \u05d1\u05de\u05e9\u05d7\u05e7 \u05d0\u05de\u05d9\u05ea\u05d9 \u05db\u05de\u05d5 \u05d0\u05d7\u05d3 \u05d4\u05d0\u05d3\u05dd. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05e8\u05d9\u05d9\u05dd \u05d4\u05d1\u05d0 \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05e4\u05e8\u05d9\u05d9\u05de\u05d9\u05dd \u05d4\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d4\u05e7\u05d5\u05d3\u05de\u05d5\u05ea \u05d5\u05d4\u05e0\u05d5\u05db\u05d7\u05d9\u05ea. \"\n\n$r6 = $d.Paragraphs.Item(6).Range\n[void]$r6.MoveEnd(1, -1)\n$r6.Text = \"\u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05de\u05ea\u05d1\u05e6\u05e2 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d3\u05d9 \u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea: \u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05de\u05e7\u05d1\u05dc \u05db\u05e7\u05dc\u05d8 \u05d0\u05ea \u05d4\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d4\u05e7\u05d5\u05d3\u05de\u05d5\u05ea \u05d0\u05d7\u05e8\u05d9 \u05d4\u05d0\u05e0\u05e7\u05d3\u05d5\u05e8 (\u05e9\u05de\u05d0\u05d5\u05de\u05df \u05d2\u05dd \u05db\u05df) \u05d5\u05d1\u05e0\u05d5\u05e1\u05e3 \u05d0\u05ea \u05d4\u05e4\u05e8\u05d9\u05d9\u05de\u05d9\u05dd \u05d4\u05e7\u05d5\u05d3\u05de\u05d9\u05dd \u05de\u05d5\u05d6\u05e0\u05d9\u05dd \u05dc\u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 (\u05d1\u05e6\u05d5\u05e8\u05d4 \u05de\u05d5\u05e8\u05e2\u05e9\u05ea \u05dc\u05e9\u05d9\u05e4\u05d5\u05e8 \u05d9\u05db\u05d5\u05dc\u05ea \u05d4\u05db\u05dc\u05dc\u05d4 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc). \u05de\u05d5\u05d3\u05dc \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05e9\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d4\u05e9\u05ea\u05de\u05e9\u05d5 \u05d1\u05d5 \u05d4\u05d9\u05e0\u05d5 \u05dc\u05d8\u05e0\u05d8\u05d9 (\u05db\u05dc\u05d5\u05de\u05e8 \u05d7\u05d9\u05d6\u05d5\u05d9 \u05d4\u05e8\u05e2\u05e9 \u05de\u05ea\u05d1\u05e6\u05e2 \u05d1\u05de\u05e8\u05d7\u05d1 \u05d4\u05dc\u05d8\u05e0\u05d8\u05d9 \u05e9\u05dc \u05d4\u05e4\u05e8\u05d9\u05d9\u05dd \u05d4\u05e0\u05d7\u05d6\u05d4). \u05e0\u05e6\u05d9\u05d9\u05df \u05db\u05d9 \u05db\u05d0\u05df \u05dc\u05d4\u05d1\u05d3\u05d9\u05dc \u05de\u05de\u05d5\u05d3\u05dc\u05d9 \u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d9\u05e9\u05e0\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d5\u05d3\u05dc \u05d4\u05d3\u05d9\u05e4\u05d5\u05d6\u05d9\u05d4 \u05d1\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05d5\u05de\u05df \u05dc\u05d7\u05d6\u05d5\u05ea \u05d0\u05ea \u05de\u05d4 \u05e9\u05e0\u05e7\u05e8\u05d0 \u05f4\u05de\u05d4\u05d9\u05e8\u05d5\u05ea\u05f4 \u05e9\u05dc \u05d4\u05e4\u05e8\u05d9\u05d9\u05dd \u05d4\u05de\u05d5\u05e8\u05e2\u05e9 \u05e9\u05d4\u05d9\u05d0 \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05e9\u05dc \u05d4\u05e4\u05e8\u05d9\u05d9\u05dd \u05d4\u05e0\u05e7\u05d9 \u05d5\u05d4\u05e8\u05e2\u05e9 \u05d4\u05de\u05ea\u05d5\u05d5\u05e1\u05e3 \u05d0\u05dc\u05d9\u05d5 \u05d1\u05d0\u05d9\u05d8\u05e8\u05e6\u05d9\u05d4. \u05e8\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d6\u05d5 \u05de\u05e9\u05e4\u05e8\u05ea \u05d0\u05ea \u05d0\u05d9\u05db\u05d5\u05ea \u05d4\u05de\u05d5\u05d3\u05dc \u05d5\u05de\u05d0\u05d9\u05e6\u05d4 \u05d4\u05ea\u05db\u05e0\u05e1\u05d5\u05ea\u05d4 (\u05de\u05d5\u05db\u05d7 \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05ea \u05db\u05e8\u05d2\u05d9\u05dc)...\"\n\n$r7 = $d.Paragraphs.Item(7).Range\n[void]$r7.MoveEnd(1, -1)\n$r7.Text = \"\u05de\u05d0\u05de\u05e8 \u05de\u05d0\u05d5\u05d3 \u05de\u05d2\u05e0\u05d9\u05d1\u2026\"\n\n$r8 = $d.Paragraphs.Item(8).Range\n[void]$r8.MoveEnd(1, -1)\n$r8.Text = \"https://arxiv.org/pdf/2408.14837 \"\n"}
